$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "contratante"
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = "Adidas"
}

[void]$ws.Range("C2:C16").Select()
